# Add data for 2024-10-29
#
# Updates the 2024 year-to-date (column K) violent-crime counts across the
# "Citywide Totals" sheet, the "By Neighborhood" summary sheet, and every
# individual neighborhood sheet that received new incidents on this date.
#
# Each entry below is the worksheet name, the column-K cell reference, and
# the new (post-update) value for that cell.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Sheet = "Citywide Totals"; Cell = "K2"; Value = 6730 },
    @{ Sheet = "Citywide Totals"; Cell = "K3"; Value = 6936 },
    @{ Sheet = "Citywide Totals"; Cell = "K4"; Value = 1436 },
    @{ Sheet = "Citywide Totals"; Cell = "K6"; Value = 7614 },
    @{ Sheet = "Citywide Totals"; Cell = "K7"; Value = 23216 },
    @{ Sheet = "By Neighborhood"; Cell = "K2"; Value = 203 },
    @{ Sheet = "By Neighborhood"; Cell = "K6"; Value = 164 },
    @{ Sheet = "By Neighborhood"; Cell = "K7"; Value = 704 },
    @{ Sheet = "By Neighborhood"; Cell = "K8"; Value = 1523 },
    @{ Sheet = "By Neighborhood"; Cell = "K9"; Value = 105 },
    @{ Sheet = "By Neighborhood"; Cell = "K11"; Value = 430 },
    @{ Sheet = "By Neighborhood"; Cell = "K15"; Value = 245 },
    @{ Sheet = "By Neighborhood"; Cell = "K19"; Value = 680 },
    @{ Sheet = "By Neighborhood"; Cell = "K20"; Value = 562 },
    @{ Sheet = "By Neighborhood"; Cell = "K21"; Value = 75 },
    @{ Sheet = "By Neighborhood"; Cell = "K22"; Value = 73 },
    @{ Sheet = "By Neighborhood"; Cell = "K27"; Value = 218 },
    @{ Sheet = "By Neighborhood"; Cell = "K29"; Value = 1262 },
    @{ Sheet = "By Neighborhood"; Cell = "K31"; Value = 255 },
    @{ Sheet = "By Neighborhood"; Cell = "K33"; Value = 1002 },
    @{ Sheet = "By Neighborhood"; Cell = "K37"; Value = 787 },
    @{ Sheet = "By Neighborhood"; Cell = "K48"; Value = 295 },
    @{ Sheet = "By Neighborhood"; Cell = "K49"; Value = 126 },
    @{ Sheet = "By Neighborhood"; Cell = "K52"; Value = 616 },
    @{ Sheet = "By Neighborhood"; Cell = "K53"; Value = 298 },
    @{ Sheet = "By Neighborhood"; Cell = "K54"; Value = 455 },
    @{ Sheet = "By Neighborhood"; Cell = "K55"; Value = 249 },
    @{ Sheet = "By Neighborhood"; Cell = "K57"; Value = 87 },
    @{ Sheet = "By Neighborhood"; Cell = "K60"; Value = 134 },
    @{ Sheet = "By Neighborhood"; Cell = "K63"; Value = 62 },
    @{ Sheet = "By Neighborhood"; Cell = "K64"; Value = 145 },
    @{ Sheet = "By Neighborhood"; Cell = "K65"; Value = 542 },
    @{ Sheet = "By Neighborhood"; Cell = "K70"; Value = 40 },
    @{ Sheet = "By Neighborhood"; Cell = "K73"; Value = 209 },
    @{ Sheet = "By Neighborhood"; Cell = "K75"; Value = 71 },
    @{ Sheet = "By Neighborhood"; Cell = "K77"; Value = 159 },
    @{ Sheet = "By Neighborhood"; Cell = "K78"; Value = 266 },
    @{ Sheet = "By Neighborhood"; Cell = "K83"; Value = 501 },
    @{ Sheet = "By Neighborhood"; Cell = "K85"; Value = 1069 },
    @{ Sheet = "By Neighborhood"; Cell = "K86"; Value = 144 },
    @{ Sheet = "By Neighborhood"; Cell = "K88"; Value = 250 },
    @{ Sheet = "By Neighborhood"; Cell = "K89"; Value = 348 },
    @{ Sheet = "By Neighborhood"; Cell = "K92"; Value = 86 },
    @{ Sheet = "By Neighborhood"; Cell = "K94"; Value = 311 },
    @{ Sheet = "By Neighborhood"; Cell = "K95"; Value = 383 },
    @{ Sheet = "By Neighborhood"; Cell = "K97"; Value = 182 },
    @{ Sheet = "By Neighborhood"; Cell = "K99"; Value = 393 },
    @{ Sheet = "By Neighborhood"; Cell = "K101"; Value = 23216 },
    @{ Sheet = "Auburn Gresham"; Cell = "K2"; Value = 233 },
    @{ Sheet = "Auburn Gresham"; Cell = "K6"; Value = 192 },
    @{ Sheet = "Auburn Gresham"; Cell = "K7"; Value = 704 },
    @{ Sheet = "Belmont Cragin"; Cell = "K2"; Value = 151 },
    @{ Sheet = "Belmont Cragin"; Cell = "K7"; Value = 430 },
    @{ Sheet = "Uptown"; Cell = "K2"; Value = 97 },
    @{ Sheet = "Uptown"; Cell = "K6"; Value = 103 },
    @{ Sheet = "Uptown"; Cell = "K7"; Value = 348 },
    @{ Sheet = "South Shore"; Cell = "K2"; Value = 352 },
    @{ Sheet = "South Shore"; Cell = "K3"; Value = 372 },
    @{ Sheet = "South Shore"; Cell = "K6"; Value = 260 },
    @{ Sheet = "South Shore"; Cell = "K7"; Value = 1069 },
    @{ Sheet = "Little Village"; Cell = "K2"; Value = 165 },
    @{ Sheet = "Little Village"; Cell = "K6"; Value = 224 },
    @{ Sheet = "Little Village"; Cell = "K7"; Value = 616 },
    @{ Sheet = "Logan Square"; Cell = "K2"; Value = 75 },
    @{ Sheet = "Logan Square"; Cell = "K7"; Value = 298 },
    @{ Sheet = "Austin"; Cell = "K2"; Value = 422 },
    @{ Sheet = "Austin"; Cell = "K3"; Value = 464 },
    @{ Sheet = "Austin"; Cell = "K6"; Value = 504 },
    @{ Sheet = "Austin"; Cell = "K7"; Value = 1523 },
    @{ Sheet = "South Chicago"; Cell = "K3"; Value = 178 },
    @{ Sheet = "South Chicago"; Cell = "K7"; Value = 501 },
    @{ Sheet = "Garfield Park"; Cell = "K3"; Value = 356 },
    @{ Sheet = "Garfield Park"; Cell = "K6"; Value = 316 },
    @{ Sheet = "Garfield Park"; Cell = "K7"; Value = 1002 },
    @{ Sheet = "West Pullman"; Cell = "K2"; Value = 132 },
    @{ Sheet = "West Pullman"; Cell = "K7"; Value = 383 },
    @{ Sheet = "Grand Crossing"; Cell = "K3"; Value = 259 },
    @{ Sheet = "Grand Crossing"; Cell = "K6"; Value = 235 },
    @{ Sheet = "Grand Crossing"; Cell = "K7"; Value = 787 },
    @{ Sheet = "New City"; Cell = "K2"; Value = 179 },
    @{ Sheet = "New City"; Cell = "K6"; Value = 198 },
    @{ Sheet = "New City"; Cell = "K7"; Value = 542 },
    @{ Sheet = "Woodlawn"; Cell = "K6"; Value = 99 },
    @{ Sheet = "Woodlawn"; Cell = "K7"; Value = 393 },
    @{ Sheet = "Gage Park"; Cell = "K2"; Value = 84 },
    @{ Sheet = "Gage Park"; Cell = "K7"; Value = 255 },
    @{ Sheet = "Lincoln Park"; Cell = "K2"; Value = 28 },
    @{ Sheet = "Lincoln Park"; Cell = "K7"; Value = 126 },
    @{ Sheet = "Loop"; Cell = "K6"; Value = 245 },
    @{ Sheet = "Loop"; Cell = "K7"; Value = 455 },
    @{ Sheet = "Englewood"; Cell = "K2"; Value = 356 },
    @{ Sheet = "Englewood"; Cell = "K3"; Value = 447 },
    @{ Sheet = "Englewood"; Cell = "K6"; Value = 368 },
    @{ Sheet = "Englewood"; Cell = "K7"; Value = 1262 },
    @{ Sheet = "Lake View"; Cell = "K2"; Value = 45 },
    @{ Sheet = "Lake View"; Cell = "K7"; Value = 295 },
    @{ Sheet = "Chatham"; Cell = "K3"; Value = 204 },
    @{ Sheet = "Chatham"; Cell = "K6"; Value = 226 },
    @{ Sheet = "Chatham"; Cell = "K7"; Value = 680 },
    @{ Sheet = "Ashburn"; Cell = "K2"; Value = 62 },
    @{ Sheet = "Ashburn"; Cell = "K7"; Value = 164 },
    @{ Sheet = "Rogers Park"; Cell = "K3"; Value = 66 },
    @{ Sheet = "Rogers Park"; Cell = "K6"; Value = 91 },
    @{ Sheet = "Rogers Park"; Cell = "K7"; Value = 266 },
    @{ Sheet = "Lower West Side"; Cell = "K3"; Value = 73 },
    @{ Sheet = "Lower West Side"; Cell = "K7"; Value = 249 },
    @{ Sheet = "Chinatown"; Cell = "K3"; Value = 19 },
    @{ Sheet = "Chinatown"; Cell = "K7"; Value = 75 },
    @{ Sheet = "Near South Side"; Cell = "K3"; Value = 40 },
    @{ Sheet = "Near South Side"; Cell = "K7"; Value = 145 },
    @{ Sheet = "Chicago Lawn"; Cell = "K2"; Value = 195 },
    @{ Sheet = "Chicago Lawn"; Cell = "K3"; Value = 181 },
    @{ Sheet = "Chicago Lawn"; Cell = "K7"; Value = 562 },
    @{ Sheet = "West Loop"; Cell = "K2"; Value = 79 },
    @{ Sheet = "West Loop"; Cell = "K6"; Value = 141 },
    @{ Sheet = "West Loop"; Cell = "K7"; Value = 311 },
    @{ Sheet = "Brighton Park"; Cell = "K6"; Value = 72 },
    @{ Sheet = "Brighton Park"; Cell = "K7"; Value = 245 },
    @{ Sheet = "Avalon Park"; Cell = "K4"; Value = 6 },
    @{ Sheet = "Avalon Park"; Cell = "K7"; Value = 105 },
    @{ Sheet = "Portage Park"; Cell = "K2"; Value = 72 },
    @{ Sheet = "Portage Park"; Cell = "K7"; Value = 209 },
    @{ Sheet = "Albany Park"; Cell = "K3"; Value = 56 },
    @{ Sheet = "Albany Park"; Cell = "K7"; Value = 203 },
    @{ Sheet = "West Town"; Cell = "K6"; Value = 97 },
    @{ Sheet = "West Town"; Cell = "K7"; Value = 182 },
    @{ Sheet = "West Elsdon"; Cell = "K2"; Value = 26 },
    @{ Sheet = "West Elsdon"; Cell = "K7"; Value = 86 },
    @{ Sheet = "O'Hare"; Cell = "K2"; Value = 19 },
    @{ Sheet = "O'Hare"; Cell = "K7"; Value = 40 },
    @{ Sheet = "United Center"; Cell = "K3"; Value = 77 },
    @{ Sheet = "United Center"; Cell = "K7"; Value = 250 },
    @{ Sheet = "Edgewater"; Cell = "K3"; Value = 51 },
    @{ Sheet = "Edgewater"; Cell = "K7"; Value = 218 },
    @{ Sheet = "Streeterville"; Cell = "K2"; Value = 25 },
    @{ Sheet = "Streeterville"; Cell = "K4"; Value = 62 },
    @{ Sheet = "Streeterville"; Cell = "K7"; Value = 144 },
    @{ Sheet = "Pullman"; Cell = "K6"; Value = 13 },
    @{ Sheet = "Pullman"; Cell = "K7"; Value = 71 },
    @{ Sheet = "Mckinley Park"; Cell = "K6"; Value = 38 },
    @{ Sheet = "Mckinley Park"; Cell = "K7"; Value = 87 },
    @{ Sheet = "Morgan Park"; Cell = "K3"; Value = 41 },
    @{ Sheet = "Morgan Park"; Cell = "K7"; Value = 134 },
    @{ Sheet = "Clearing"; Cell = "K4"; Value = 4 },
    @{ Sheet = "Clearing"; Cell = "K7"; Value = 73 },
    @{ Sheet = "Riverdale"; Cell = "K2"; Value = 67 },
    @{ Sheet = "Riverdale"; Cell = "K7"; Value = 159 }
)

foreach ($update in $updates) {
    $ws = $wb.Worksheets.Item($update.Sheet)
    $ws.Range($update.Cell).Value = $update.Value
}
